$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 168787
$ws.Range("C4").Value = 159631
$ws.Range("C5").Value = 9157
$ws.Range("C7").Value = 5.43
$ws.Range("C8").Value = 65.45
